# Apply the ValueSet-fr-medication-reconciliation-type.xlsx updates described in the diff:
#  - URL (ValueSet) changed from .../fhir/fr/medication/... to .../ig/fhir/medication/...
#  - Name changed from FrMedicationReconciliationType to FRMedicationReconciliationType
#  - Title casing fix: InterOp'Santé -> Interop'Santé
#  - Date bumped to 2026-01-15T08:54:26+00:00
#  - Jurisdiction value filled in as FRANCE
#  - System URI (CodeSystem) changed from .../fhir/fr/medication/... to .../ig/fhir/medication/...

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value  = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-type"
$wsMeta.Range("B4").Value  = "FRMedicationReconciliationType"
$wsMeta.Range("B5").Value  = "value set Interop'Santé - Type d'écart/erreur sur une ligne de traitement d'une FCT"
$wsMeta.Range("B8").Value  = "2026-01-15T08:54:26+00:00"
$wsMeta.Range("B11").Value = "FRANCE"

$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-type"
